# Daily_Motivations.xlsx - append new daily baseline/score rows for 2025-02-10
# (mirrors the existing per-day block of sleep / activity / weekly_activity rows)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-DateTextCell($addr, $text) {
    # Force text storage so a date-looking string (e.g. "2025-02-10") is not
    # auto-converted into a date serial number, matching the existing rows
    # which store the date as plain text. Resetting the style back to
    # "Normal" afterwards keeps the cell free of any explicit style index,
    # just like the pre-existing data rows.
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $text
    $ws.Range($addr).Style = "Normal"
}

$newDate = "2025-02-10"

# Row 29: sleep
Set-DateTextCell "A29" $newDate
$ws.Range("B29").Value = "sleep"
$ws.Range("C29").Value = $false
$ws.Range("D29").Value = $false

# Row 30: activity
Set-DateTextCell "A30" $newDate
$ws.Range("B30").Value = "activity"
$ws.Range("C30").Value = $true
$ws.Range("D30").Value = $true

# Row 31: weekly_activity
Set-DateTextCell "A31" $newDate
$ws.Range("B31").Value = "weekly_activity"
$ws.Range("C31").Value = $false
$ws.Range("D31").Value = $false
